# Fruta / hortaliza, semanal
# Insert a new weekly record at row 40 (pushing the existing rows 40-42
# down to 41-43) and populate it with the latest price report.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 40:42 down to 41:43, duplicating the formatting of row 40
# (this also extends the date-format style onto the new D40 cell).
$ws.Rows.Item(40).Insert()

$ws.Cells.Item(40, 1).Value = 1
$ws.Cells.Item(40, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(40, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(40, 4).Value = 45142
$ws.Cells.Item(40, 5).Value = 15
$ws.Cells.Item(40, 6).Value = 100112013
$ws.Cells.Item(40, 7).Value = "Alcachofa"
$ws.Cells.Item(40, 8).Value = "Madrigal"
$ws.Cells.Item(40, 9).Value = "Primera"
$ws.Cells.Item(40, 10).Value = 170
$ws.Cells.Item(40, 11).Value = 20000
$ws.Cells.Item(40, 12).Value = 22000
$ws.Cells.Item(40, 13).Value = 21176
$ws.Cells.Item(40, 14).Value = "$/caja 40 unidades"
$ws.Cells.Item(40, 15).Value = "Región de Coquimbo"
$ws.Cells.Item(40, 16).Value = 529
$ws.Cells.Item(40, 17).Value = 40
$ws.Cells.Item(40, 18).Value = "Hortaliza"
